$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.809.97"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "2.105.70"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.03"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.96"
$ws.Range("E7").Value = "  +2.02%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.98"
$ws.Range("E12").Value = "  +6.86%  "

$ws.Range("D13").Value = "2.417.84"
$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.02"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.803"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").Value = "2.081.43"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "38.890.89"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.66"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "0.0₃0844"
$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  -5.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.79"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("E28").Value = "  -0.56%  "

$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.36"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("E31").Value = "  +9.50%  "

$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.58"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.17"
$ws.Range("E34").Value = "  +11.76%  "

$ws.Range("E35").Value = "  +1.10%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.05"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +3.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.11"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").Value = "1.529.02"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("E44").Value = "  +7.36%  "

$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0917"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("E48").Value = "  +5.17%  "

$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").Value = "2.305.01"
$ws.Range("E51").Value = "  +0.66%  "
